$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 1048
$ws.Range("F7").Value = 2696
$ws.Range("F9").Value = 1324
$ws.Range("F13").Value = 1204
$ws.Range("F19").Value = 545
$ws.Range("F22").Value = 664
$ws.Range("F24").Value = 233
$ws.Range("F28").Value = 626
$ws.Range("F29").Value = 6850
$ws.Range("F34").Value = 191
$ws.Range("F37").Value = 116
$ws.Range("F42").Value = 22
$ws.Range("F45").Value = 151
$ws.Range("F47").Value = 126

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 57

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 757

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 757
$ws.Range("F6").Value = 1048
$ws.Range("F7").Value = 2696
$ws.Range("F9").Value = 1324
$ws.Range("F13").Value = 1204
$ws.Range("F21").Value = 545
$ws.Range("F24").Value = 57
$ws.Range("F25").Value = 664
$ws.Range("F27").Value = 233
$ws.Range("F30").Value = 626
$ws.Range("F31").Value = 6850
$ws.Range("F36").Value = 191
$ws.Range("F44").Value = 22
$ws.Range("F46").Value = 151
$ws.Range("F48").Value = 126
